$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-14 down to 10-15
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the new market data entry
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "Vega Monumental Concepción"
$ws.Range("C9").Value = "Bíobío"
$ws.Range("D9").Value = 44484
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100107
$ws.Range("H9").Value = "Otros"
$ws.Range("I9").Value = 100107002
$ws.Range("J9").Value = "Chirimoya"
$ws.Range("K9").Value = "Cultivar IV Región"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 25000
$ws.Range("O9").Value = 26000
$ws.Range("P9").Value = 25500
$ws.Range("Q9").Value = "$/bandeja 10 kilos"
$ws.Range("R9").Value = "Provincia de Limarí"
$ws.Range("S9").Value = 2550
$ws.Range("T9").Value = 10

# Ensure D9 keeps the same date number format as the other date cells (s="2")
$ws.Range("D9").NumberFormat = $ws.Range("D10").NumberFormat
